$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("13-09-2021", 312.2, 381, 196, 305, 84.59999999999999, 113.1, 13.7, 153.6, 468, 1483, 291.4, 278, 140, 348, 174),
    @("14-09-2021", 315.5, 384, 200, 308, 83.90000000000001, 115.4, 13, 156.4, 472.9, 1501.1, 295.1, 281, 145, 347, 178),
    @("15-09-2021", 312.3, 380, 199, 305, 83.90000000000001, 116.9, 16.8, 152.1, 473.2, 1489.8, 291.4, 277, 143, 341, 174),
    @("16-09-2021", 312.2, 381, 198, 306, 83.09999999999999, 116.5, 12.8, 148.6, 481.2, 1495.8, 290.3, 279, 142, 344, 172),
    @("17-09-2021", 311.8, 381, 197, 307, 83.7, 117.4, 9.199999999999999, 147.9, 486.9, 1512.2, 288.4, 279, 141, 343, 169)
)

$startRow = 176
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
